# Update column G ("K") values for rows 2-18 in the active sheet.
# These are the "Strike#"/"K" stat values recomputed per the commit message
# ("regen save_data to use K instead of Strike#, regen std/mean, calc and
# write s_vals").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 4
    3  = 1
    4  = 1
    5  = 1
    6  = 2
    7  = 0
    8  = 0
    9  = 2
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 1
    16 = 1
    17 = 0
    18 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
